$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07490066666666667
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 0.03072872253977778
$ws.Range("R2").Value = 0.276558502858
$ws.Range("S2").Value = 0.003499619873322347
$ws.Range("T2").Value = 0.003499619873322347

# Row 3
$ws.Range("G3").Value = 0.07490066666666667
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("R3").Value = 68.69016980833801
$ws.Range("S3").Value = 0.8692174743460166
$ws.Range("T3").Value = 0.8692174743460165

# Row 4
$ws.Range("G4").Value = 0.07490066666666667
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("S4").Value = 0.1272829057806611
$ws.Range("T4").Value = 0.1272829057806611

$wb.Save()
